$d = $word.ActiveDocument

# The paragraph that currently ends the document ("... // 배성환 대리").
$lastPara = $d.Paragraphs.Last
$fullRange = $lastPara.Range

# Replacement WordprocessingML: the same paragraph but with its paragraph-mark
# run-formatting (<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> inside <w:pPr>)
# stripped, followed by a page break and the new list content.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00EB1855" w:rsidRDefault="00EB1855" w:rsidP="00EB1855"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>1</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">월 </w:t></w:r><w:r><w:t>16</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">일엔 이 설정파일을 백업하는방법 적용하는방법 진행예정 </w:t></w:r><w:r><w:t>&gt;&gt; TFTP(UDP)</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>방식 사용한다구함</w:t></w:r><w:r w:rsidR="0028666C"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0028666C"><w:t xml:space="preserve">// </w:t></w:r><w:r w:rsidR="0028666C"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>배성환 대리</w:t></w:r></w:p>
<w:p><w:pPr><w:widowControl/><w:wordWrap/><w:autoSpaceDE/><w:autoSpaceDN/></w:pPr><w:r><w:br w:type="page"/></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">TX: Transmitt ( </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">송신 </w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:t xml:space="preserve">RX: receive ( </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">수신 </w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>광케이블</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">싱글모드 </w:t></w:r><w:r><w:t xml:space="preserve">( </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>멀리나감</w:t></w:r><w:r><w:t xml:space="preserve"> ) </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">멀티모드 </w:t></w:r><w:r><w:t xml:space="preserve">( </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">짧게나감 </w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>모듈</w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>컨버터)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>G</w:t></w:r><w:r><w:t xml:space="preserve">BIC ( </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">큰거 </w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>S</w:t></w:r><w:r><w:t>FP (</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 작은거 </w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>G</w:t></w:r><w:r><w:t xml:space="preserve">BIC </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">일 경우 </w:t></w:r><w:r><w:t xml:space="preserve">SC </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>타입</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>S</w:t></w:r><w:r><w:t>FP</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 일 경우 L</w:t></w:r><w:r><w:t xml:space="preserve">C </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>타입</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$fullRange.InsertXML($xml) | Out-Null

# InsertXML on a range that reaches the document's final paragraph mark can't
# actually remove that mark, so Word leaves a blank trailing paragraph behind
# (carrying the original "배성환 대리" paragraph's pPr/rPr). Merge it away so
# the new "...LC 타입" paragraph (with its own bookmark) is the true last one.
$n = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($n - 1)
$mergeRange = $d.Range($secondLast.Range.End - 1, $d.Content.End)
$mergeRange.Delete() | Out-Null
